# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets, reflecting the newer scrape snapshot referenced in the
# commit message ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 614
$ws1.Range("F6").Value = 2812
$ws1.Range("F14").Value = 5923
$ws1.Range("F18").Value = 103
$ws1.Range("F23").Value = 21

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F28").Value = 130
$ws2.Range("F32").Value = 193

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2591
$ws3.Range("F12").Value = 630

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2591
$ws4.Range("F12").Value = 614
$ws4.Range("F13").Value = 2812
$ws4.Range("F16").Value = 630
$ws4.Range("F22").Value = 5923
$ws4.Range("F25").Value = 103
$ws4.Range("F39").Value = 130
$ws4.Range("F45").Value = 193
